# Generate Report for Handback
# - Flips the "Ready for handoff" status to "Handed back: in sync with en-US"
#   on every sheet that surfaces it (Overview + each language sheet).
# - Fills in the "Latest Target File" / "Latest Handback File" /
#   "Latest Handback DateTime" columns for each language sheet now that the
#   handback has happened, including the hyperlinks for the two file columns.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: columns B (zh-cn) and C (de-de) show the same status
# text for both rows.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

# ---------------------------------------------------------------------
# Per-language sheets.
# ---------------------------------------------------------------------
$languages = @(
  @{
    Sheet = "zh-cn"
    HandbackDateTime = "2016-03-19 14:50:28"
    Row2 = @{
      SourceDisplay = "06128c65-284a-4f9a-af02-d4ace1ef9822.md"
      SourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/233cdb63adf21d7be95332c90ea0a8324ec00108/e2e/06128c65-284a-4f9a-af02-d4ace1ef9822.md"
      TargetDisplay = "06128c65-284a-4f9a-af02-d4ace1ef9822.6bb1de5248d32704318862486bb0c69045d9137d.zh-cn.xlf"
      TargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/465b3028c155af283c03e2a76224575c13f93641/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/06128c65-284a-4f9a-af02-d4ace1ef9822.6bb1de5248d32704318862486bb0c69045d9137d.zh-cn.xlf"
    }
    Row3 = @{
      SourceDisplay = "7b11a892-6e38-4768-8f00-0d924dfa8f6f.md"
      SourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/233cdb63adf21d7be95332c90ea0a8324ec00108/e2e/7b11a892-6e38-4768-8f00-0d924dfa8f6f.md"
      TargetDisplay = "7b11a892-6e38-4768-8f00-0d924dfa8f6f.9c32c9c5e3b59884e68a64650affd715de2476c4.zh-cn.xlf"
      TargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/465b3028c155af283c03e2a76224575c13f93641/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/7b11a892-6e38-4768-8f00-0d924dfa8f6f.9c32c9c5e3b59884e68a64650affd715de2476c4.zh-cn.xlf"
    }
  },
  @{
    Sheet = "de-de"
    HandbackDateTime = "2016-03-19 14:50:33"
    Row2 = @{
      SourceDisplay = "06128c65-284a-4f9a-af02-d4ace1ef9822.md"
      SourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/233cdb63adf21d7be95332c90ea0a8324ec00108/e2e/06128c65-284a-4f9a-af02-d4ace1ef9822.md"
      TargetDisplay = "06128c65-284a-4f9a-af02-d4ace1ef9822.6bb1de5248d32704318862486bb0c69045d9137d.de-de.xlf"
      TargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fc425b92c81b1fc426b738598c96e265ac8676a6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/06128c65-284a-4f9a-af02-d4ace1ef9822.6bb1de5248d32704318862486bb0c69045d9137d.de-de.xlf"
    }
    Row3 = @{
      SourceDisplay = "7b11a892-6e38-4768-8f00-0d924dfa8f6f.md"
      SourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/233cdb63adf21d7be95332c90ea0a8324ec00108/e2e/7b11a892-6e38-4768-8f00-0d924dfa8f6f.md"
      TargetDisplay = "7b11a892-6e38-4768-8f00-0d924dfa8f6f.9c32c9c5e3b59884e68a64650affd715de2476c4.de-de.xlf"
      TargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fc425b92c81b1fc426b738598c96e265ac8676a6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/7b11a892-6e38-4768-8f00-0d924dfa8f6f.9c32c9c5e3b59884e68a64650affd715de2476c4.de-de.xlf"
    }
  }
)

# BGR packed color that round-trips to the workbook's existing hyperlink
# font color (RGB FF6495ED).
$hyperlinkColor = 15570276

foreach ($lang in $languages) {
  $ws = $wb.Worksheets.Item($lang.Sheet)

  # Status column (C) -> handed back.
  $ws.Range("C2").Value = $statusNew
  $ws.Range("C3").Value = $statusNew

  # Latest Handback DateTime (H) -> real timestamp instead of the zero date.
  $ws.Range("H2").Value = $lang.HandbackDateTime
  $ws.Range("H3").Value = $lang.HandbackDateTime

  # Latest Target File (F) + Latest Handback File (G) for row 2.
  $ws.Range("F2").Value = $lang.Row2.SourceDisplay
  $ws.Hyperlinks.Add($ws.Range("F2"), $lang.Row2.SourceUrl, "", "", $lang.Row2.SourceDisplay)
  $ws.Range("F2").Font.Underline = 2
  $ws.Range("F2").Font.Color = $hyperlinkColor

  $ws.Range("G2").Value = $lang.Row2.TargetDisplay
  $ws.Hyperlinks.Add($ws.Range("G2"), $lang.Row2.TargetUrl, "", "", $lang.Row2.TargetDisplay)
  $ws.Range("G2").Font.Underline = 2
  $ws.Range("G2").Font.Color = $hyperlinkColor

  # Latest Target File (F) + Latest Handback File (G) for row 3.
  $ws.Range("F3").Value = $lang.Row3.SourceDisplay
  $ws.Hyperlinks.Add($ws.Range("F3"), $lang.Row3.SourceUrl, "", "", $lang.Row3.SourceDisplay)
  $ws.Range("F3").Font.Underline = 2
  $ws.Range("F3").Font.Color = $hyperlinkColor

  $ws.Range("G3").Value = $lang.Row3.TargetDisplay
  $ws.Hyperlinks.Add($ws.Range("G3"), $lang.Row3.TargetUrl, "", "", $lang.Row3.TargetDisplay)
  $ws.Range("G3").Font.Underline = 2
  $ws.Range("G3").Font.Color = $hyperlinkColor
}
